# Reorder the "Recorded By" (column G) values so that the email address
# appears before the literal "System"/"system" tokens.
#
# Only these exact strings are affected (verified against the target diff):
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
#   "system, System, backup@backdoor.com" -> "backup@backdoor.com, system, System"
#
# Other combinations (e.g. "System, admin@admin.com", lone "System", or
# rows that already have the email first) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "system, System, backup@backdoor.com" = "backup@backdoor.com, system, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    $value = $cell.Value2
    if ($null -ne $value -and $map.ContainsKey($value)) {
        $cell.Value = $map[$value]
    }
}
